$d = $word.ActiveDocument

# Locate the paragraph that ends with the drone pathfinding sentence, then
# insert a brand-new paragraph right after it (before the trailing empty
# paragraph), matching the formatting (sz 24 / szCs 24) used throughout this
# section of the report.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*pathfinding towards the player without going through walls anymore.*") {
        $anchor = $p
        break
    }
}

$anchor.Range.InsertParagraphAfter()

$newPara = $anchor.Next()
$aposRight = [char]0x2019
$newPara.Range.Text = "In addition to this I also had issues when setting up the networking this is because it is something that I am new to, this meant that wrapping my head around it wasn" + $aposRight + "t very easy. The issues I had were when sending high scores to be saved on the server instead of retrieving the high scores from the server."
